$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0.2250718607779158
$ws.Range("B2").Value = -0.3042742479111624
$ws.Range("A3").Value = -0.4387633823875903
$ws.Range("B3").Value = -0.3666699275062994
$ws.Range("A4").Value = -0.5009261959971308
$ws.Range("B4").Value = -0.4474214250777439
$ws.Range("A5").Value = -0.1530744787944337
$ws.Range("B5").Value = -0.2200429065265428
$ws.Range("A6").Value = -0.1657178361506949
$ws.Range("B6").Value = -0.1443424340210184
$ws.Range("A7").Value = -0.2963440572473632
$ws.Range("B7").Value = -0.3006785005963426
$ws.Range("A8").Value = -0.5278915093624672
$ws.Range("B8").Value = -0.4167639099711564
$ws.Range("A9").Value = -0.5627393231283633
$ws.Range("B9").Value = -0.4622878142841325
$ws.Range("A10").Value = -0.4968082027901717
$ws.Range("B10").Value = -0.2803444740328301
$ws.Range("A11").Value = -0.2914590922241845
$ws.Range("B11").Value = -0.298744234038695
$ws.Range("A12").Value = -0.2314086778818056
$ws.Range("B12").Value = -0.1990670213938235
$ws.Range("A13").Value = -0.04773324398888493
$ws.Range("B13").Value = -0.07393797950775674
$ws.Range("A14").Value = -0.2939329335286706
$ws.Range("B14").Value = -0.2169967485155734
$ws.Range("A15").Value = -0.1430035065739936
$ws.Range("B15").Value = -0.04406473220843669
$ws.Range("A16").Value = -0.1893388278177911
$ws.Range("B16").Value = -0.07991983670123821
$ws.Range("A17").Value = 0.05951402412171154
$ws.Range("B17").Value = 0.1387323640026155
$ws.Range("A18").Value = 0.03708600989497735
$ws.Range("B18").Value = 0.05239589056464296
$ws.Range("A19").Value = 0.06069684514083785
$ws.Range("B19").Value = 0.1007067326617474
$ws.Range("A20").Value = -0.1572175032763125
$ws.Range("B20").Value = -0.08198769250113755
$ws.Range("A21").Value = -0.003509973189657789
$ws.Range("B21").Value = 0.0382361029445268
$ws.Range("A22").Value = 0.03161691759604291
$ws.Range("B22").Value = 0.1086412327525411
$ws.Range("A23").Value = -0.006455082552080424
$ws.Range("B23").Value = 0.02083511503494362
$ws.Range("A24").Value = 0.8134419455026002
$ws.Range("B24").Value = 0.6275414214716868
$ws.Range("A25").Value = 0.1250648350631644
$ws.Range("B25").Value = 0.09369898545125437
$ws.Range("A26").Value = 0.1411182489151249
$ws.Range("B26").Value = 0.1277892020732319
$ws.Range("A27").Value = 0.09688272276242904
$ws.Range("B27").Value = 0.1094231370355884
$ws.Range("A28").Value = 0.3134806690642843
$ws.Range("B28").Value = 0.1994971955143899
$ws.Range("A29").Value = 0.6787123820577274
$ws.Range("B29").Value = 0.5500951161456252
$ws.Range("A30").Value = 0.2093280834551174
$ws.Range("B30").Value = 0.1693179181437547
$ws.Range("A31").Value = 0.01527400927096381
$ws.Range("B31").Value = 0.03363750468044523
$ws.Range("A32").Value = 0.181117137605032
$ws.Range("B32").Value = 0.1816233813767025
$ws.Range("A33").Value = 0.1181998612117385
$ws.Range("B33").Value = 0.130918521061543
$ws.Range("A34").Value = 0.08767057010953486
$ws.Range("B34").Value = 0.05295205450232404
$ws.Range("A35").Value = 0.4354695729439043
$ws.Range("B35").Value = 0.2888409977143332
$ws.Range("A36").Value = 0.2912862454769364
$ws.Range("B36").Value = 0.1533495673690221
$ws.Range("A37").Value = 0.05039469239145366
$ws.Range("B37").Value = -0.0064213426474793
$ws.Range("A38").Value = 0.3449913081183461
$ws.Range("B38").Value = 0.2993136825837875
$ws.Range("A39").Value = -0.08949323932445991
$ws.Range("B39").Value = -0.1543204714276162
$ws.Range("A40").Value = 0.1644700705108344
$ws.Range("B40").Value = 0.1810125953524771
$ws.Range("A41").Value = -0.1659332337896862
$ws.Range("B41").Value = -0.1951428843044119
$ws.Range("A42").Value = 0.2550643703135704
$ws.Range("B42").Value = 0.2531039654461633
$ws.Range("A43").Value = 0.1670444740067412
$ws.Range("B43").Value = 0.1673087483328268
$ws.Range("A44").Value = -0.1037455053460379
$ws.Range("B44").Value = -0.08646426511031464
$ws.Range("A45").Value = -0.104794556577724
$ws.Range("B45").Value = -0.08798903293491311
$ws.Range("A46").Value = -0.1925971413069253
$ws.Range("B46").Value = -0.1761233106356346
$ws.Range("A47").Value = -0.1900851859062858
$ws.Range("B47").Value = -0.1758569071786635
$ws.Range("A48").Value = -0.2281306172181959
$ws.Range("B48").Value = -0.2028629821515599
$ws.Range("A49").Value = -0.2203072234534432
$ws.Range("B49").Value = -0.2004981598542138
$ws.Range("A50").Value = -0.1494051867186948
$ws.Range("B50").Value = -0.1368666296624268
$ws.Range("A51").Value = -0.2130587950665102
$ws.Range("B51").Value = -0.2213222295619814
$ws.Range("A52").Value = -0.2130587950665102
$ws.Range("B52").Value = -0.2213222295619814
$ws.Range("A53").Value = -0.2040477838392285
$ws.Range("B53").Value = -0.1859288231682639
$ws.Range("A54").Value = -0.2095204110480183
$ws.Range("B54").Value = -0.1994429908687647
$ws.Range("A55").Value = -0.1720922102463877
$ws.Range("B55").Value = -0.1490445105651763
$ws.Range("A56").Value = -0.1587711423644652
$ws.Range("B56").Value = -0.1409908660108635
$ws.Range("A57").Value = -0.1972692006844608
$ws.Range("B57").Value = -0.155633558866462
$ws.Range("A58").Value = -0.1813738256061986
$ws.Range("B58").Value = -0.213516391298937
$ws.Range("A59").Value = -0.2222151474455981
$ws.Range("B59").Value = -0.2355574678922932
$ws.Range("A60").Value = -0.2095247580315729
$ws.Range("B60").Value = -0.2389769316631077
$ws.Range("A61").Value = -0.2722358043226074
$ws.Range("B61").Value = -0.229451432604436
$ws.Range("A62").Value = -0.1606000546365942
$ws.Range("B62").Value = -0.1007039242574877
$ws.Range("A63").Value = -0.3529740683774172
$ws.Range("B63").Value = -0.4020249240780268
$ws.Range("A64").Value = -0.2577263942527296
$ws.Range("B64").Value = -0.2638823367836596
$ws.Range("A65").Value = -0.1998195694491312
$ws.Range("B65").Value = -0.2303599148437729
$ws.Range("A66").Value = -0.120372322293627
$ws.Range("B66").Value = -0.09740285784665913
$ws.Range("A67").Value = -0.06217892092423713
$ws.Range("B67").Value = -0.05586642611577743

Write-Host "done"